# Fruta / hortaliza, semanal
#
# Two new weekly price rows are inserted in the "Palta" data block, just
# above the former row 594. All subsequent rows (594-655) shift down by
# two, becoming 596-657, and the sheet's used range grows from A1:T655 to
# A1:T657.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 594 (pushes old 594.. down to 596..)
$ws.Rows.Item(594).Insert()
$ws.Rows.Item(594).Insert()

# --- New row 594 ---
$ws.Range("A594").Value = 7
$ws.Range("B594").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C594").Value = "Ñuble"
$ws.Range("D594").Value = 44918
$ws.Range("E594").Value = 16
$ws.Range("F594").Value = "Fruta"
$ws.Range("G594").Value = 100106
$ws.Range("H594").Value = "Oleaginosos"
$ws.Range("I594").Value = 100106002
$ws.Range("J594").Value = "Palta"
$ws.Range("K594").Value = "Hass"
$ws.Range("L594").Value = "Primera"
$ws.Range("M594").Value = 100
$ws.Range("N594").Value = 3400
$ws.Range("O594").Value = 3500
$ws.Range("P594").Value = 3450
$ws.Range("Q594").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R594").Value = "Provincia de Quillota"
$ws.Range("S594").Value = 3450
$ws.Range("T594").Value = 1

# --- New row 595 ---
$ws.Range("A595").Value = 7
$ws.Range("B595").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C595").Value = "Ñuble"
$ws.Range("D595").Value = 44918
$ws.Range("E595").Value = 16
$ws.Range("F595").Value = "Fruta"
$ws.Range("G595").Value = 100106
$ws.Range("H595").Value = "Oleaginosos"
$ws.Range("I595").Value = 100106002
$ws.Range("J595").Value = "Palta"
$ws.Range("K595").Value = "Hass"
$ws.Range("L595").Value = "Segunda"
$ws.Range("M595").Value = 50
$ws.Range("N595").Value = 3000
$ws.Range("O595").Value = 3000
$ws.Range("P595").Value = 3000
$ws.Range("Q595").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R595").Value = "Provincia de Quillota"
$ws.Range("S595").Value = 3000
$ws.Range("T595").Value = 1

# Make sure the date cells keep the workbook's date number format (style
# index 2 == numFmtId 165, "YYYY-MM-DD HH:MM:SS") — Insert() already
# carries the format down from the row above, but set it explicitly too.
$ws.Range("D594:D595").NumberFormat = "YYYY-MM-DD HH:MM:SS"
